$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (comb1)
$ws.Range("B2").Value = 1
$ws.Range("E2").Value = 1.7

# Update row 3 (comb2)
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 0

# Update row 4 (comb3)
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.405

# Add new row 5 (comb4) - copy the formatted label cell A4 into A5 to
# carry over the same style, then set the new values
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "comb4"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0.35
